$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shared strings shift: "is_locked"/"is_enabled" columns removed; "order_by"
# and "rem" move left into B1/C1; a new "tenant_id" label takes D1; the old
# E1 cell (which held the old "rem" string) is removed entirely.
$ws.Cells.Item(1, 2).Value = '<%=comment.order_by%>'
$ws.Cells.Item(1, 3).Value = '<%=comment.rem%>'
$ws.Cells.Item(1, 4).Value = '<%=comment.tenant_id_lbl%><%selectList.tenant_id = data.findAllTenant.map((item) => item.lbl)%><%_dataValidation_({ sqref: `${ _col }2:${ _col }${ _lastRow }`, formula1: `"${ selectList.tenant_id.join(",") }"` })%>'

$ws.Cells.Item(1, 5).ClearContents()
